$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 7 to the "Logs" sheet ---
$ws.Range("A7").Value = "Waarom zit er verschil tussen de EcoPro-600 en EcoPro-700?"
$ws.Range("B7").Value = "mailmind.test@zohomail.eu"
$ws.Range("C7").Value = "Testmail #6: Waarom zit er verschil tussen de EcoPro-600 en EcoPro-700?"
$ws.Range("D7").Value = "Productinformatie"
$ws.Range("E7").Value = @"
Beste klant,
Bedankt voor uw vraag over de EcoPro-600 en EcoPro-700. Het verschil tussen deze twee modellen zit voornamelijk in de capaciteit en functionaliteiten. De EcoPro-700 heeft bijvoorbeeld een grotere watertank en een extra reinigingsfunctie ten opzichte van de EcoPro-600.
Indien u meer specifieke informatie wilt over de verschillen tussen deze modellen, kunt u de productpagina’s op onze website raadplegen of contact met ons opnemen voor een gedetailleerdere uitleg.
Met vriendelijke groet,
[Naam] 
E-mailassistent 
--------------------------------------------------------------------------
  Dit is een testmail. Gelieve hier niet op te antwoorden.
"@
$ws.Range("F7").Value = "2025-06-30 19:49:18"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Ja"
$ws.Range("J7").Value = "Nee"

# --- Extend the existing conditional formatting ranges so they include row 7,
#     preserving rule order, priorities and dxf (style) references ---
$ws.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D7"))
$ws.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G7"))
$ws.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H7"))
$ws.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I7"))
$ws.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J7"))

# --- Update the "Dashboard" summary sheet so the rows are resorted by count (desc) ---
$dashboard.Range("A2").Value = "Productinformatie"
$dashboard.Range("B2").Value = 2
$dashboard.Range("A3").Value = "Retour / Terugbetaling"
$dashboard.Range("B3").Value = 2
$dashboard.Range("A4").Value = "Bestelling / Levering"
$dashboard.Range("B4").Value = 1
$dashboard.Range("A5").Value = "Openingstijden / Locatie"
$dashboard.Range("B5").Value = 1
